# Apply crypto price/volume updates to match the target diff (Wed Nov 29 2023 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.188.48"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "2.060.09"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'229.83"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("D7").Value = "'59.73"
$ws.Range("E7").Value = "  +8.01%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  +3.22%  "
$ws.Range("D10").Value = "'0.0813"
$ws.Range("E10").Value = "  +4.34%  "
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").Value = "'14.76"
$ws.Range("E12").Value = "  +5.15%  "
$ws.Range("D13").Value = "2.362.53"
$ws.Range("E13").Value = "  +2.35%  "
$ws.Range("D14").Value = "'21.24"
$ws.Range("E14").Value = "  +7.50%  "
$ws.Range("D15").Value = "'0.758"
$ws.Range("E15").Value = "  +2.82%  "
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "2.055.99"
$ws.Range("E17").Value = "  +1.90%  "
$ws.Range("D18").Value = "38.069.38"
$ws.Range("E18").Value = "  +2.70%  "
$ws.Range("D19").Value = "'6.29"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").Value = "'69.97"
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").Value = "0.0₃0836"
$ws.Range("E21").Value = "  +2.97%  "
$ws.Range("D22").Value = "'225.35"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").Value = "'2.44"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("E25").Value = "  +4.00%  "
$ws.Range("D26").Value = "'166.42"
$ws.Range("E26").Value = "  +1.26%  "
$ws.Range("D27").Value = "'9.27"
$ws.Range("E27").Value = "  +3.77%  "
$ws.Range("D28").Value = "'0.133"
$ws.Range("E28").Value = "  +7.26%  "
$ws.Range("D29").Value = "'19.07"
$ws.Range("E29").Value = "  +2.58%  "
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("E31").Value = "  +2.85%  "
$ws.Range("D32").Value = "'4.57"
$ws.Range("E32").Value = "  +3.69%  "
$ws.Range("D33").Value = "'4.60"
$ws.Range("E33").Value = "  +3.36%  "
$ws.Range("E34").Value = "  +10.57%  "
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("E36").Value = "  +0.84%  "
$ws.Range("D37").Value = "'6.19"
$ws.Range("E37").Value = "  +15.17%  "
$ws.Range("D38").Value = "'3.29"
$ws.Range("E38").Value = "  +5.45%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "1.534.35"
$ws.Range("E40").Value = "  +4.91%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'17.19"
$ws.Range("E41").Value = "  +7.54%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'98.37"
$ws.Range("E42").Value = "  +3.72%  "
$ws.Range("E43").Value = "  +2.43%  "
$ws.Range("E44").Value = "  +4.37%  "
$ws.Range("D45").Value = "'0.0925"
$ws.Range("E45").Value = "  +1.58%  "
$ws.Range("E46").Value = "  +1.47%  "
$ws.Range("D47").Value = "'4.12"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("D50").Value = "'7.13"
$ws.Range("D51").Value = "2.250.68"
$ws.Range("E51").Value = "  +2.47%  "
